# Auto-generated script to apply F-column ('想去人数') value updates
# across the 4 worksheets, per the target diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
# Expected sheet name: 展览
$ws.Range("F4").Value = 276
$ws.Range("F6").Value = 2029
$ws.Range("F8").Value = 2056
$ws.Range("F10").Value = 1276
$ws.Range("F11").Value = 2959
$ws.Range("F12").Value = 1817
$ws.Range("F13").Value = 1248
$ws.Range("F15").Value = 60
$ws.Range("F16").Value = 0
$ws.Range("F17").Value = 592
$ws.Range("F18").Value = 1201
$ws.Range("F19").Value = 1922
$ws.Range("F20").Value = 1201
$ws.Range("F21").Value = 184
$ws.Range("F22").Value = 2172
$ws.Range("F24").Value = 646
$ws.Range("F25").Value = 4716
$ws.Range("F26").Value = 50
$ws.Range("F30").Value = 76
$ws.Range("F32").Value = 421
$ws.Range("F33").Value = 82
$ws.Range("F34").Value = 0
$ws.Range("F38").Value = 39
$ws.Range("F39").Value = 3408
$ws.Range("F41").Value = 1102
$ws.Range("F45").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("F48").Value = 756
$ws.Range("F49").Value = 34

$ws = $wb.Worksheets.Item(2)
# Expected sheet name: 演出
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 190
$ws.Range("F7").Value = 148544
$ws.Range("F8").Value = 148544
$ws.Range("F15").Value = 383
$ws.Range("F18").Value = 3
$ws.Range("F19").Value = 254
$ws.Range("F30").Value = 73
$ws.Range("F31").Value = 311
$ws.Range("F32").Value = 14
$ws.Range("F33").Value = 0
$ws.Range("F35").Value = 35
$ws.Range("F38").Value = 118
$ws.Range("F42").Value = 17

$ws = $wb.Worksheets.Item(3)
# Expected sheet name: 本地生活
$ws.Range("F4").Value = 3213
$ws.Range("F10").Value = 297
$ws.Range("F12").Value = 819

$ws = $wb.Worksheets.Item(4)
# Expected sheet name: 全部类型
$ws.Range("F5").Value = 297
$ws.Range("F9").Value = 2029
$ws.Range("F10").Value = 148544
$ws.Range("F13").Value = 1276
$ws.Range("F14").Value = 2959
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 1716
$ws.Range("F19").Value = 145
$ws.Range("F20").Value = 3
$ws.Range("F21").Value = 592
$ws.Range("F23").Value = 1201
$ws.Range("F24").Value = 1922
$ws.Range("F25").Value = 1201
$ws.Range("F27").Value = 184
$ws.Range("F32").Value = 0
$ws.Range("F35").Value = 137
$ws.Range("F40").Value = 35
$ws.Range("F41").Value = 421
$ws.Range("F45").Value = 3408
$ws.Range("F48").Value = 1102
$ws.Range("F49").Value = 0
$ws.Range("F50").Value = 0
$ws.Range("F51").Value = 756
$ws.Range("F52").Value = 48
$ws.Range("F53").Value = 4
